$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header correction ---------------------------------------------
# B2 incorrectly referenced the stray "unnamed: 1_level_1" label; it should
# read "total" just like B1.
$ws.Range("B2").Value = "total"

# --- Remove the empty "situação do domicílio" / "grandes regiões" header
#     rows and shift the real data up so every label lines up with its row.
# Row 5 ("situação do domicílio", no data) -> becomes "urbana" + its data
$ws.Range("A5").Value = "urbana"
$ws.Range("B5").Value2 = 0
$ws.Range("C5").Value2 = 2.3
$ws.Range("D5").Value2 = 7.79
$ws.Range("E5").Value2 = 4.58
$ws.Range("F5").Value2 = 5.99
$ws.Range("G5").Value2 = 2.43

# Row 6 ("urbana") -> becomes "rural" + its data
$ws.Range("A6").Value = "rural"
$ws.Range("B6").Value2 = 0
$ws.Range("C6").Value2 = 8.24
$ws.Range("D6").Value2 = 26.07
$ws.Range("E6").Value2 = 11.86
$ws.Range("F6").Value2 = 14.68
$ws.Range("G6").Value2 = 6.21

# Row 7 ("rural") -> becomes "norte" + its data
$ws.Range("A7").Value = "norte"
$ws.Range("B7").Value2 = 0
$ws.Range("C7").Value2 = 4.76
$ws.Range("D7").Value2 = 19.23
$ws.Range("E7").Value2 = 9.77
$ws.Range("F7").Value2 = 14.13
$ws.Range("G7").Value2 = 5.88

# Row 8 ("grandes regiões", no data) -> becomes "nordeste" + its data
$ws.Range("A8").Value = "nordeste"
$ws.Range("B8").Value2 = 0
$ws.Range("C8").Value2 = 3.94
$ws.Range("D8").Value2 = 13.48
$ws.Range("E8").Value2 = 7.26
$ws.Range("F8").Value2 = 9.59
$ws.Range("G8").Value2 = 3.73

# Row 9 ("norte") -> becomes "sudeste" + its data
$ws.Range("A9").Value = "sudeste"
$ws.Range("B9").Value2 = 0
$ws.Range("C9").Value2 = 4.26
$ws.Range("D9").Value2 = 12.58
$ws.Range("E9").Value2 = 8.32
$ws.Range("F9").Value2 = 10.41
$ws.Range("G9").Value2 = 4.29

# Row 10 ("nordeste") -> becomes "sul" + its data
$ws.Range("A10").Value = "sul"
$ws.Range("B10").Value2 = 0
$ws.Range("C10").Value2 = 5.82
$ws.Range("D10").Value2 = 20
$ws.Range("E10").Value2 = 11.48
$ws.Range("F10").Value2 = 14.38
$ws.Range("G10").Value2 = 5.98

# Rows 11-12 ("sudeste"/"sul") held the data that has now moved up into
# rows 9-10, so the trailing rows are deleted entirely (sheet shrinks from
# A1:G12 to A1:G10).
$ws.Range("A11:G12").Delete()
